$d = $word.ActiveDocument

# The footer's first paragraph currently reads:
#   "Copyright 2018, " / "Created " / "by Babak Aryan"
# It needs to become a page-number field (inserted via the "Page Numbers
# (Bottom of Page)" building block) followed by two tabs and the text
# "Created by Babak Aryan," (with "Babak" flagged by the spell checker).
$footer = $d.Sections(1).Footers(1)
$para1 = $footer.Range.Paragraphs(1)
$rng = $para1.Range.Duplicate
[void]$rng.MoveEnd(1, -1)   # exclude the paragraph mark so the <w:p> itself
                             # (and its paraId/rsid identity) is left
                             # untouched - only its run content is replaced

$newPara1 = '<w:p>' +
  '<w:sdt>' +
    '<w:sdtPr>' +
      '<w:rPr><w:noProof/></w:rPr>' +
      '<w:id w:val="-820351100"/>' +
      '<w:docPartObj>' +
        '<w:docPartGallery w:val="Page Numbers (Bottom of Page)"/>' +
        '<w:docPartUnique/>' +
      '</w:docPartObj>' +
    '</w:sdtPr>' +
    '<w:sdtContent>' +
      '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
      '<w:r><w:instrText xml:space="preserve"> PAGE   \* MERGEFORMAT </w:instrText></w:r>' +
      '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
      '<w:r><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r>' +
      '<w:r><w:rPr><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' +
    '</w:sdtContent>' +
  '</w:sdt>' +
  '<w:r><w:rPr><w:noProof/></w:rPr><w:tab/></w:r>' +
  '<w:r><w:rPr><w:noProof/></w:rPr><w:tab/></w:r>' +
  '<w:r><w:t xml:space="preserve">Created by </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Babak</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> Aryan,</w:t></w:r>' +
  '</w:p>'

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' + $newPara1 + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$rng.InsertXML($xmlFrag)

# The second footer paragraph, "For Wintriss Technical Schools", becomes a
# copyright/license notice.
[void]$footer.Range.Find.Execute("For Wintriss Technical Schools", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "Copyright 2018 under the terms of a Creative Commons License", 2)
